$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the introductory paragraph: "medewerker" -> "stagiair"
$ws.Range("C3").Value = "Ben je bezig met het aannemen van een nieuwe stagiair? Deze checklist helpt je om het onboardingproces overzichtelijk en effectief te laten verlopen. Het zorgt ervoor dat alle belangrijke stappen – van voorbereiding tot de eerste werkdag – goed zijn geregeld."

# Update the checklist item: "werkdag" -> "stagedag"
$ws.Range("C12").Value = "Duidelijkheid over eerste stagedag"

# Update the view (scroll position / selection) to match the saved state
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("C15").Select()
